# Apply corrected values produced by the new version of readxl.
# Each entry maps a cell reference to its corrected value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "P2" = 686; "X2" = 345; "Y2" = 3208; "Z2" = 1604; "AA2" = "406 x 178 x 54"; "AC2" = 4899; "AD2" = 2481; "AE2" = 1.41; "AG2" = 0.42; "AH2" = 2035; "AJ2" = 4899; "AK2" = 345000; "AL2" = 45763; "AP2" = 4757; "AR2" = 2450; "AS2" = 21141; "AT2" = 0.34; "AV2" = 0.97; "AW2" = 2371; "AY2" = 45763; "AZ2" = 42; "BA2" = 806
    "P3" = 686; "X3" = 286; "Y3" = 2984; "Z3" = 1492; "AA3" = "356 x 171 x 45"; "AC3" = 4068; "AD3" = 3096; "AE3" = 1.15; "AG3" = 0.5600000000000001; "AH3" = 2298; "AJ3" = 4068; "AK3" = 286500; "AL3" = 73309; "AP3" = 4036; "AR3" = 2034; "AS3" = 16809; "AT3" = 0.35; "AV3" = 0.97; "AW3" = 1965; "AY3" = 73309; "AZ3" = 27; "BA3" = 749
    "P4" = 1164; "Y4" = 9712; "Z4" = 4856; "AA4" = "686 x 254 x 140"; "AB4" = 345; "AC4" = 12282; "AD4" = 18040; "AE4" = 0.83; "AG4" = 0.78; "AH4" = 9593; "AI4" = 345; "AJ4" = 12282; "AK4" = 890000; "AL4" = 118056; "AM4" = 0.32; "AP4" = 11940; "AQ4" = 345; "AR4" = 6141; "AS4" = 107362; "AT4" = 0.24; "AV4" = 0.99; "AW4" = 6088; "AY4" = 118056; "AZ4" = 127; "BA4" = 2441
    "P5" = 1164; "Y5" = 9156; "Z5" = 4578; "AA5" = "610 x 229 x 125"; "AB5" = 345; "AC5" = 10971; "AD5" = 25230; "AE5" = 0.66; "AG5" = 0.87; "AH5" = 9501; "AI5" = 345; "AJ5" = 10971; "AK5" = 795000; "AL5" = 203424; "AM5" = 0.23; "AP5" = 10893; "AQ5" = 345; "AR5" = 5486; "AS5" = 81454; "AT5" = 0.26; "AV5" = 0.99; "AW5" = 5413; "AY5" = 203424; "AZ5" = 84; "BA5" = 2297
    "P6" = 1962; "Y6" = 16368; "Z6" = 8184; "AA6" = "914 x 305 x 224"; "AC6" = 19734; "AD6" = 49875; "AE6" = 0.63; "AG6" = 0.88; "AH6" = 17339; "AJ6" = 19734; "AK6" = 1430000; "AL6" = 189686; "AP6" = 19185; "AR6" = 9867; "AS6" = 232133; "AT6" = 0.21; "AV6" = 1; "AW6" = 9854; "AY6" = 189686; "AZ6" = 214; "BA6" = 4113
    "P7" = 1962; "Y7" = 15432; "Z7" = 7716; "AA7" = "1016 x 305 x 222"; "AC7" = 19527; "AD7" = 104398; "AE7" = 0.43; "AG7" = 0.9399999999999999; "AH7" = 18434; "AJ7" = 19527; "AK7" = 1415000; "AL7" = 362068; "AP7" = 19388; "AR7" = 9764; "AS7" = 197935; "AT7" = 0.22; "AV7" = 1; "AW7" = 9716; "AY7" = 362068; "AZ7" = 142; "BA7" = 3872
    "P8" = 1724; "Y8" = 14384; "Z8" = 7192; "AA8" = "762 x 267 x 197"; "AC8" = 17319; "AD8" = 31835; "AE8" = 0.74; "AG8" = 0.83; "AH8" = 14364; "AJ8" = 17319; "AK8" = 1255000; "AL8" = 166473; "AP8" = 16837; "AR8" = 8660; "AS8" = 169333; "AT8" = 0.23; "AW8" = 8610; "AY8" = 166473; "AZ8" = 188; "BA8" = 3615
    "P9" = 1724; "Y9" = 13560; "Z9" = 6780; "AA9" = "610 x 305 x 179"; "AC9" = 15732; "AD9" = 39149; "AE9" = 0.63; "AG9" = 0.88; "AH9" = 13792; "AJ9" = 15732; "AK9" = 1140000; "AL9" = 291702; "AP9" = 15620; "AR9" = 7866; "AS9" = 236278; "AT9" = 0.18; "AV9" = 1; "AW9" = 7896; "AY9" = 291702; "AZ9" = 125; "BA9" = 3402
    "P10" = 960; "Y10" = 8010; "Z10" = 4005; "AA10" = "686 x 254 x 125"; "AC10" = 10971; "AD10" = 15652; "AE10" = 0.84; "AG10" = 0.77; "AH10" = 8489; "AJ10" = 10971; "AK10" = 795000; "AL10" = 105455; "AP10" = 10666; "AR10" = 5486; "AS10" = 90781; "AT10" = 0.25; "AW10" = 5430; "AY10" = 105455; "AZ10" = 104; "BA10" = 2013
    "P11" = 960; "Y11" = 7552; "Z11" = 3776; "AA11" = "533 x 210 x 101"; "AC11" = 8901; "AD11" = 15737; "AE11" = 0.75; "AG11" = 0.82; "AH11" = 7316; "AJ11" = 8901; "AK11" = 645000; "AL11" = 165042; "AP11" = 8838; "AR11" = 4450; "AS11" = 55753; "AT11" = 0.28; "AV11" = 0.98; "AW11" = 4368; "AY11" = 165042; "AZ11" = 70; "BA11" = 1895
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
